$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H ("divorces") values added for this data draft
$hValues = @{
    8 = 93
    14 = 211
    20 = 66
    26 = 54
    32 = 87
    38 = 31
    44 = 316
    50 = 227
    56 = 13
    62 = 258
    68 = 35
    74 = 120
    80 = 38
    86 = 94
    92 = 68
    98 = 137
    104 = 58
    110 = 169
    116 = 15
    122 = 35
    128 = 1106
    134 = 41
    140 = 11
    146 = 131
    152 = 34
    158 = 62
    164 = 90
    170 = 38
    176 = 3559
    182 = 183
    188 = 50
    194 = 54
    200 = 215
    206 = 359
    212 = 103
    218 = 173
    224 = 52
    230 = 91
    236 = 48
    242 = 46
    248 = 136
    254 = 343
    260 = 1001
    266 = 206
    272 = 5
    278 = 93
    284 = 174
    290 = 18
    296 = 71
    302 = 46
    308 = 23
    314 = 108
    320 = 655
    326 = 77
    332 = 44
    338 = 23
    344 = 26
    350 = 101
    356 = 30
    362 = 283
    368 = 40
    374 = 58
    380 = 36
    386 = 18
    392 = 41
    398 = 72
    404 = 81
    410 = 68
    416 = 288
    422 = 70
    428 = 67
    434 = 117
    440 = 860
    446 = 172
    452 = 22
    458 = 13
    464 = 304
    470 = 85
    476 = 438
    482 = 170
    488 = 91
    494 = 150
    500 = 170
    506 = 83
    512 = 321
    518 = 380
    524 = 20
    530 = 21
    536 = 229
    542 = 294
    548 = 279
    554 = 21
    560 = 250
    566 = 155
    572 = 198
}

foreach ($row in $hValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $hValues[$row]
}

# H342 previously used the one-off highlight style (fillId 2); restore it to
# the plain "applyFill" style shared by the rest of the sheet now that the
# highlight is gone.
$ws.Range("H342").Style = "Normal"
